# Update Shift Unit design — fill in activity-log rows 7-13 on the
# "Activity Log - Part 2" sheet with dates, start/end times and
# descriptions, then leave the sheet scrolled/selected at G13 (matches
# the author's last interactive selection before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log - Part 2")

# Dates / times first (column order doesn't affect shared-string order).
$ws.Range("C7").Value2  = 43931
$ws.Range("D7").Value2  = 0.83333333333333337
$ws.Range("E7").Value2  = 1

$ws.Range("C8").Value2  = 43932
$ws.Range("D8").Value2  = 0.41666666666666669
$ws.Range("E8").Value2  = 0.5

$ws.Range("C9").Value2  = 43932
$ws.Range("D9").Value2  = 0.54166666666666663
$ws.Range("E9").Value2  = 0.66666666666666663

$ws.Range("C10").Value2 = 43933
$ws.Range("D10").Value2 = 0.375
$ws.Range("E10").Value2 = 0.45833333333333331

$ws.Range("C11").Value2 = 43933
$ws.Range("D11").Value2 = 0.83333333333333337
$ws.Range("E11").Value2 = 0.91666666666666663

$ws.Range("C12").Value2 = 43934
$ws.Range("D12").Value2 = 0.83333333333333337
$ws.Range("E12").Value2 = 0.89583333333333337

$ws.Range("C13").Value2 = 43935
$ws.Range("D13").Value2 = 0.75
$ws.Range("E13").Value2 = 0.83333333333333337

# Descriptions (column G) — written in the order the author actually
# typed them so new shared-string entries land at the same indices.
$ws.Range("G7").Value  = "Designed Execution Unit "
$ws.Range("G9").Value  = "Designed Shift Unit"
$ws.Range("G10").Value = "Debug Shift Unit"
$ws.Range("G8").Value  = "Debug Execution Unit "
$ws.Range("G11").Value = "Created circuit diagrams for Shift Unit and Execution Unit"
$ws.Range("G12").Value = "Shift Unit Documentation (Sections: Functional Behaviour, VHDL Interface)"
$ws.Range("G13").Value = "Updating design of Shift Unit"

$ws.Activate()
$ws.Range("G13").Select()
